# Add a new paragraph "Emery" right after the existing last paragraph ("Matt"),
# mirroring the formatting of the rest of the document.

$d = $word.ActiveDocument

# Locate the last paragraph in the document (the one ending in "Matt").
$lastPara = $d.Paragraphs($d.Paragraphs.Count)

# Insert a brand-new paragraph mark right after it.
$lastPara.Range.InsertParagraphAfter()

# The newly created paragraph is now the last paragraph; give it the text "Emery".
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.InsertAfter("Emery")
